$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.416.84"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.851.62"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "241.05"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "0.6297"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D8").Value = "0.07685"
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("D9").Value = "0.2938"
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").Value = "0.07751"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "1.853.36"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "0.00001093"
$ws.Range("D15").Value = "0.6813"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "83.69"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "2.102.53"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "6.161"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").Value = "29.460.11"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "229.41"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D23").Value = "7.458"
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "157.11"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "0.1389"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").Value = "8.390"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").Value = "17.69"
$ws.Range("D29").Value = "1.315"
$ws.Range("E29").Value = "  +4.05%  "
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").Value = "0.05722"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("D33").Value = "4.051"
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").Value = "1.849"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").Value = "1.164"
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("D36").Value = "0.7086"
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "2.777"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D40").Value = "1.219.53"
$ws.Range("E40").Value = "  -2.43%  "
$ws.Range("D41").Value = "6.492"
$ws.Range("E41").Value = "  +5.04%  "
$ws.Range("D42").Value = "0.9081"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.011.43"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "101.66"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "66.40"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.00000000120"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "7.136"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value = "0.4018"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "8.971"
$ws.Range("E51").Value = "  -1.22%  "
